# Update test-data paths (commit: "Update paths for test data.")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages_with_footer")

# Row 3: espanol node path -> new espanol press-release path
$ws.Range("A3").Value = "espanol/noticias/comunicados-de-prensa/2018/oropharyngeal-hpv-cisplatin"

# Row 5: node path -> annual report slug
$ws.Range("A5").Value = "annual-report-nation-2018"

# Row 10 was a duplicate data row; clear it out entirely
$ws.Range("A10:C10").ClearContents()

# Header row (path / type / language) is now bold
$ws.Range("A1:C1").Font.Bold = $true

